# Update RAD FEIN/SSN object-recognition test-data timestamps, and move the
# "currently active" tab/selection from Existing -> NewTaxReturn, matching a
# newer Katalon test run (Oct 31 2023) that replaced the Oct 9-11 2023 run.

$wb = $excel.ActiveWorkbook

# --- Existing sheet: B2:B12 -----------------------------------------------
$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value  = "Tue Oct 31 17:52:13 EDT 2023"
$ws.Range("B3").Value  = "Tue Oct 31 17:52:49 EDT 2023"
$ws.Range("B4").Value  = "Tue Oct 31 17:53:25 EDT 2023"
$ws.Range("B5").Value  = "Tue Oct 31 17:54:01 EDT 2023"
$ws.Range("B6").Value  = "Tue Oct 31 17:54:37 EDT 2023"
$ws.Range("B7").Value  = "Tue Oct 31 17:55:13 EDT 2023"
$ws.Range("B8").Value  = "Tue Oct 31 17:55:51 EDT 2023"
$ws.Range("B9").Value  = "Tue Oct 31 17:56:29 EDT 2023"
$ws.Range("B10").Value = "Tue Oct 31 17:57:07 EDT 2023"
$ws.Range("B11").Value = "Tue Oct 31 17:57:44 EDT 2023"
$ws.Range("B12").Value = "Tue Oct 31 17:58:21 EDT 2023"

# --- Extension sheet: B2:B7 -------------------------------------------------
$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Tue Oct 31 17:58:58 EDT 2023"
$ws.Range("B3").Value = "Tue Oct 31 17:59:33 EDT 2023"
$ws.Range("B4").Value = "Tue Oct 31 18:00:07 EDT 2023"
$ws.Range("B5").Value = "Tue Oct 31 18:00:41 EDT 2023"
$ws.Range("B6").Value = "Tue Oct 31 18:01:15 EDT 2023"
$ws.Range("B7").Value = "Tue Oct 31 18:01:50 EDT 2023"

# --- NewTaxReturn sheet: B2:B16 --------------------------------------------
$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value  = "Tue Oct 31 18:02:25 EDT 2023"
$ws.Range("B3").Value  = "Tue Oct 31 18:03:00 EDT 2023"
$ws.Range("B4").Value  = "Tue Oct 31 18:03:34 EDT 2023"
$ws.Range("B5").Value  = "Tue Oct 31 18:04:09 EDT 2023"
$ws.Range("B6").Value  = "Tue Oct 31 18:04:43 EDT 2023"
$ws.Range("B7").Value  = "Tue Oct 31 18:05:18 EDT 2023"
$ws.Range("B8").Value  = "Tue Oct 31 18:05:52 EDT 2023"
$ws.Range("B9").Value  = "Tue Oct 31 18:06:27 EDT 2023"
$ws.Range("B10").Value = "Tue Oct 31 18:07:01 EDT 2023"
$ws.Range("B11").Value = "Tue Oct 31 18:07:36 EDT 2023"
$ws.Range("B12").Value = "Tue Oct 31 18:08:10 EDT 2023"
$ws.Range("B13").Value = "Tue Oct 31 18:08:44 EDT 2023"
$ws.Range("B14").Value = "Tue Oct 31 18:09:19 EDT 2023"
$ws.Range("B15").Value = "Tue Oct 31 18:09:53 EDT 2023"
$ws.Range("B16").Value = "Tue Oct 31 18:10:28 EDT 2023"

# Column B on NewTaxReturn got wider to fit the longer recognition label.
$ws.Columns.Item(2).ColumnWidth = 29.5

# --- Personal_IND sheet: B2:B6 ----------------------------------------------
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Tue Oct 31 18:11:37 EDT 2023"
$ws.Range("B3").Value = "Tue Oct 31 18:12:08 EDT 2023"
$ws.Range("B4").Value = "Tue Oct 31 18:12:40 EDT 2023"
$ws.Range("B5").Value = "Tue Oct 31 18:13:12 EDT 2023"
$ws.Range("B6").Value = "Tue Oct 31 18:13:44 EDT 2023"

# --- Personal_EL sheet: B2 --------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Tue Oct 31 18:11:02 EDT 2023"

# --- Tab / selection state --------------------------------------------------
# Existing is no longer the active tab; its selection moves to the whole
# updated Date column.
$wsExisting = $wb.Worksheets.Item("Existing")
$wsExisting.Range("C2:C12").Select()

# NewTaxReturn becomes the active tab, selection on its Date column.
$wsNewTaxReturn = $wb.Worksheets.Item("NewTaxReturn")
$wsNewTaxReturn.Select()
$wsNewTaxReturn.Range("C2:C16").Select()
